$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value-only corrections in column C (re-synced counts) ---
$ws.Cells.Item(905, 3).Value = 20
$ws.Cells.Item(917, 3).Value = 18
$ws.Cells.Item(931, 3).Value = 17
$ws.Cells.Item(951, 3).Value = 31
$ws.Cells.Item(956, 3).Value = 36
$ws.Cells.Item(959, 3).Value = 7
$ws.Cells.Item(961, 3).Value = 30
$ws.Cells.Item(974, 3).Value = 12
$ws.Cells.Item(975, 3).Value = 35
$ws.Cells.Item(979, 3).Value = 17
$ws.Cells.Item(980, 3).Value = 36
$ws.Cells.Item(982, 3).Value = 24
$ws.Cells.Item(997, 3).Value = 15
$ws.Cells.Item(998, 3).Value = 38
$ws.Cells.Item(1001, 3).Value = 31
$ws.Cells.Item(1002, 3).Value = 28
$ws.Cells.Item(1007, 3).Value = 12
$ws.Cells.Item(1008, 3).Value = 28
$ws.Cells.Item(1019, 3).Value = 7
$ws.Cells.Item(1025, 3).Value = 16
$ws.Cells.Item(1030, 3).Value = 19
$ws.Cells.Item(1031, 3).Value = 36
$ws.Cells.Item(1036, 3).Value = 17
$ws.Cells.Item(1037, 3).Value = 36

# --- Rewrite rows 1038-1055: data re-grouped/shifted, two brand new rows added at the end (1054,1055) for new date 44173 ---
$ws.Cells.Item(1038, 1).Value = 44170
$ws.Cells.Item(1038, 2).Value = "50-59"
$ws.Cells.Item(1038, 3).Value = 3
$ws.Cells.Item(1039, 1).Value = 44170
$ws.Cells.Item(1039, 2).Value = "60-69"
$ws.Cells.Item(1039, 3).Value = 11
$ws.Cells.Item(1040, 1).Value = 44170
$ws.Cells.Item(1040, 2).Value = "70-79"
$ws.Cells.Item(1040, 3).Value = 15
$ws.Cells.Item(1041, 1).Value = 44170
$ws.Cells.Item(1041, 2).Value = "80+"
$ws.Cells.Item(1041, 3).Value = 35
$ws.Cells.Item(1042, 1).Value = 44171
$ws.Cells.Item(1042, 2).Value = "50-59"
$ws.Cells.Item(1042, 3).Value = 1
$ws.Cells.Item(1043, 1).Value = 44171
$ws.Cells.Item(1043, 2).Value = "60-69"
$ws.Cells.Item(1043, 3).Value = 9
$ws.Cells.Item(1044, 1).Value = 44171
$ws.Cells.Item(1044, 2).Value = "70-79"
$ws.Cells.Item(1044, 3).Value = 18
$ws.Cells.Item(1045, 1).Value = 44171
$ws.Cells.Item(1045, 2).Value = "80+"
$ws.Cells.Item(1045, 3).Value = 22
$ws.Cells.Item(1046, 1).Value = 44172
$ws.Cells.Item(1046, 2).Value = "30-39"
$ws.Cells.Item(1046, 3).Value = 1
$ws.Cells.Item(1047, 1).Value = 44172
$ws.Cells.Item(1047, 2).Value = "40-49"
$ws.Cells.Item(1047, 3).Value = 1
$ws.Cells.Item(1048, 1).Value = 44172
$ws.Cells.Item(1048, 2).Value = "50-59"
$ws.Cells.Item(1048, 3).Value = 3
$ws.Cells.Item(1049, 1).Value = 44172
$ws.Cells.Item(1049, 2).Value = "60-69"
$ws.Cells.Item(1049, 3).Value = 11
$ws.Cells.Item(1050, 1).Value = 44172
$ws.Cells.Item(1050, 2).Value = "70-79"
$ws.Cells.Item(1050, 3).Value = 17
$ws.Cells.Item(1051, 1).Value = 44172
$ws.Cells.Item(1051, 2).Value = "80+"
$ws.Cells.Item(1051, 3).Value = 32
$ws.Cells.Item(1052, 1).Value = 44173
$ws.Cells.Item(1052, 2).Value = "50-59"
$ws.Cells.Item(1052, 3).Value = 1
$ws.Cells.Item(1053, 1).Value = 44173
$ws.Cells.Item(1053, 2).Value = "60-69"
$ws.Cells.Item(1053, 3).Value = 5
$ws.Cells.Item(1054, 1).Value = 44173
$ws.Cells.Item(1054, 2).Value = "70-79"
$ws.Cells.Item(1054, 3).Value = 3
$ws.Cells.Item(1055, 1).Value = 44173
$ws.Cells.Item(1055, 2).Value = "80+"
$ws.Cells.Item(1055, 3).Value = 17

# --- Ensure new date cells (rows that previously had no data) carry the same date number format as the rest of column A ---
$dateFormat = $ws.Cells.Item(1037, 1).NumberFormat
$ws.Cells.Item(1054, 1).NumberFormat = $dateFormat
$ws.Cells.Item(1055, 1).NumberFormat = $dateFormat

Write-Host "Applied COVID deaths-by-age dataset refresh; sheet now spans A1:C1055"
